$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column so numeric-looking strings are not
# auto-converted to Number type by Excel (matches original inline-string cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "93.177.90"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "3.433.63"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "231.21"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "620.24"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -2.79%  "
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "0.962"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "3.432.22"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "42.55"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "4.068.57"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "93.002.03"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "8.14"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "3.429.71"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "18.11"
$ws.Range("E20").Value = "  +4.21%  "
$ws.Range("D21").Value = "11.61"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "501.52"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "3.35"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").Value = "0.443"
$ws.Range("E24").Value = "  -7.32%  "
$ws.Range("D25").Value = "6.59"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("E26").Value = "  -3.78%  "
$ws.Range("D27").Value = "91.91"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "3.610.67"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "11.93"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "11.42"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("D37").Value = "0.547"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "556.48"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.40"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "0.923"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "23.67"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "3.68"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "5.50"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0409"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "53.68"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "2.11"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").Value = "8.03"
$ws.Range("E51").Value = "  +0.09%  "
